$d = $word.ActiveDocument

# --- Tail: merge the bookmark-only paragraph into the "Prototipos" paragraph,
# and drop the trailing empty paragraph, reverting to the pre-review layout.
$pProto = $d.Paragraphs.Item(14)
$pBookmark = $d.Paragraphs.Item(15)

# Remove the " " text content of the bookmark-only paragraph (keep its mark,
# which anchors the bookmarkStart/bookmarkEnd, for now).
$spaceRange = $d.Range($pBookmark.Range.Start, $pBookmark.Range.End - 1)
if (($spaceRange.End - $spaceRange.Start) -gt 0) {
    $spaceRange.Delete()
}

# Merge the (now empty) bookmark paragraph into the "Prototipos" paragraph by
# deleting the "Prototipos" paragraph's own trailing mark.
$pProto = $d.Paragraphs.Item(14)
$markRange = $d.Range($pProto.Range.End - 1, $pProto.Range.End)
$markRange.Delete()

# Drop the final empty paragraph that used to close the document.
$pProto = $d.Paragraphs.Item(14)
$trailingMark = $d.Range($pProto.Range.End - 1, $pProto.Range.End)
$trailingMark.Delete()

# --- Remove the three reviewer comment paragraphs (and the blank spacer
# paragraphs that followed two of them), restoring the original text.

# "LAS CARACTERISTAS NOS PARECE QUE ESTAN BIEN ES LO QUE  ESPECIFICAMENTE  NOS VA A PEDIR EL APLICATIVO"
$d.Paragraphs.Item(13).Range.Delete()

# blank spacer paragraph after "NOS PARECE QUE LAS NECESIDADES ESTAN BIEN"
$d.Paragraphs.Item(11).Range.Delete()
# "NOS PARECE QUE LAS NECESIDADES ESTAN BIEN"
$d.Paragraphs.Item(10).Range.Delete()

# blank spacer paragraph after "Debido a que la aplicación ... entendible."
$d.Paragraphs.Item(8).Range.Delete()
# "Debido a que la aplicación de las pausas activas ... entendible."
$d.Paragraphs.Item(7).Range.Delete()
